# JU_Total_Subtotal.xlsx - 1.0.3
# Replace the old ":For:staffs.$size:staffs.$get(index).XXX" loop-template
# strings in the header/template row (row 3) with the new shorthand
# ":staffs[].XXX" collection-expansion syntax, and move the saved
# cursor/selection from C10 to C8.
#
# Edited cells (row 3 of the single worksheet):
#   D3 : staffName  template
#   E3 : time.YMD   template
#   C3 : staffNo    template
# (B3 ":orgName", F3 ":staffs.$size" etc. are untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = ":staffs[].staffName"
$ws.Range("E3").Value = ":staffs[].time.YMD"
$ws.Range("C3").Value = ":staffs[].staffNo"

# Move the saved selection/active cell to C8 (was C10).
$ws.Range("C8").Select() | Out-Null
